$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete rows (3,4,5); only the header + one data row remain.
$ws.Rows("3:5").Delete() | Out-Null

# Row 2 gets refreshed with the new TPM-derived values.
$ws.Range("A2").Value2 = "FAPs"
$ws.Range("B2").Value2 = "Cxcl5"
$ws.Range("C2").Value2 = "Cxcr1"
$ws.Range("D2").Value2 = "FAPs"

$ws.Range("E2").Value2 = 3
$ws.Range("F2").Value2 = 1
$ws.Range("G2").Value2 = 4.799114
$ws.Range("H2").Value2 = 14.397342
$ws.Range("I2").Value2 = 1
$ws.Range("J2").Value2 = 1
$ws.Range("K2").Value2 = 1
$ws.Range("L2").Value2 = 0.3333333333333333
$ws.Range("M2").Value2 = 0.001809666666666667
$ws.Range("N2").Value2 = 0.005429
$ws.Range("O2").Value2 = 1
$ws.Range("P2").Value2 = 1
$ws.Range("Q2").Value2 = 0.008684796635333334
$ws.Range("R2").Value2 = 0.078163169718
$ws.Range("S2").Value2 = 1
$ws.Range("T2").Value2 = 1
